$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mdk"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06163433333333333
$ws.Range("N2").Value = 0.184903
$ws.Range("O2").Value = 0.006690894379667537
$ws.Range("P2").Value = 0.006690894379667537
$ws.Range("Q2").Value = 0.03685570829588888
$ws.Range("R2").Value = 0.3317013746629999
$ws.Range("S2").Value = 0.0002236535351136539
$ws.Range("T2").Value = 0.0002236535351136539

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mdk"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1030763333333333
$ws.Range("N3").Value = 0.309229
$ws.Range("O3").Value = 0.01118975126488057
$ws.Range("P3").Value = 0.01118975126488057
$ws.Range("Q3").Value = 0.06163693298988888
$ws.Range("R3").Value = 0.5547323969089999
$ws.Range("S3").Value = 0.0003740348129000616
$ws.Range("T3").Value = 0.0003740348129000616

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mdk"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.046962666666667
$ws.Range("N4").Value = 27.140888
$ws.Range("O4").Value = 0.9821193543554519
$ws.Range("P4").Value = 0.9821193543554518
$ws.Range("Q4").Value = 5.409845437983111
$ws.Range("R4").Value = 48.688608941848
$ws.Range("S4").Value = 0.03282886457939433
$ws.Range("T4").Value = 0.03282886457939433

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mdk"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06163433333333333
$ws.Range("N5").Value = 0.184903
$ws.Range("O5").Value = 0.006690894379667537
$ws.Range("P5").Value = 0.006690894379667537
$ws.Range("Q5").Value = 0.9784538964658888
$ws.Range("R5").Value = 8.806085068192999
$ws.Range("S5").Value = 0.005937605950574971
$ws.Range("T5").Value = 0.005937605950574971

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mdk"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1030763333333333
$ws.Range("N6").Value = 0.309229
$ws.Range("O6").Value = 0.01118975126488057
$ws.Range("P6").Value = 0.01118975126488057
$ws.Range("Q6").Value = 1.636351600299889
$ws.Range("R6").Value = 14.727164402699
$ws.Range("S6").Value = 0.009929963010282946
$ws.Range("T6").Value = 0.009929963010282946

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mdk"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.046962666666667
$ws.Range("N7").Value = 27.140888
$ws.Range("O7").Value = 0.9821193543554519
$ws.Range("P7").Value = 0.9821193543554518
$ws.Range("Q7").Value = 143.6218320803031
$ws.Range("R7").Value = 1292.596488722728
$ws.Range("S7").Value = 0.8715483150229518
$ws.Range("T7").Value = 0.8715483150229517

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Mdk"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.06163433333333333
$ws.Range("N8").Value = 0.184903
$ws.Range("O8").Value = 0.006690894379667537
$ws.Range("P8").Value = 0.006690894379667537
$ws.Range("Q8").Value = 0.08727816059733333
$ws.Range("R8").Value = 0.785503445376
$ws.Range("S8").Value = 0.0005296348939789121
$ws.Range("T8").Value = 0.0005296348939789121

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Mdk"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1030763333333333
$ws.Range("N9").Value = 0.309229
$ws.Range("O9").Value = 0.01118975126488057
$ws.Range("P9").Value = 0.01118975126488057
$ws.Range("Q9").Value = 0.1459626848853333
$ws.Range("R9").Value = 1.313664163968
$ws.Range("S9").Value = 0.0008857534416975658
$ws.Range("T9").Value = 0.0008857534416975658

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Mdk"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.046962666666667
$ws.Range("N10").Value = 27.140888
$ws.Range("O10").Value = 0.9821193543554519
$ws.Range("P10").Value = 0.9821193543554518
$ws.Range("Q10").Value = 12.81107814161067
$ws.Range("R10").Value = 115.299703274496
$ws.Range("S10").Value = 0.07774217475310584
$ws.Range("T10").Value = 0.07774217475310584
